$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-11) for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# These values correspond to a reshuffle of the existing rows' data.
$rowsData = @{
    2  = @{ D = 44547; J = 200; K = 13000; L = 14000; M = 13500; P = 750 }
    3  = @{ D = 45005; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    4  = @{ D = 44568; J = 500; K = 15000; L = 16000; M = 15500; P = 861 }
    5  = @{ D = 44964; J = 300; K = 20000; L = 21000; M = 20500; P = 1139 }
    6  = @{ D = 44977; J = 400; K = 16500; L = 17000; M = 16750; P = 931 }
    7  = @{ D = 44984; J = 200; K = 17000; L = 18000; M = 17500; P = 972 }
    8  = @{ D = 44998; J = 320; K = 17000; L = 18000; M = 17500; P = 972 }
    9  = @{ D = 44557; J = 400; K = 13000; L = 14000; M = 13500; P = 750 }
    10 = @{ D = 44960; J = 400; K = 19500; L = 20000; M = 19750; P = 1097 }
    11 = @{ D = 44957; J = 400; K = 21000; L = 22000; M = 21500; P = 1194 }
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 10).Value = $vals.J
    $ws.Cells.Item($r, 11).Value = $vals.K
    $ws.Cells.Item($r, 12).Value = $vals.L
    $ws.Cells.Item($r, 13).Value = $vals.M
    $ws.Cells.Item($r, 16).Value = $vals.P
}
